$p = $ppt.ActivePresentation

# NOTE: TextRange.Text normalises curly quotes/apostrophes to their plain
# ASCII equivalents when *read*, even though the underlying XML (and the
# rendered slide) keeps the curly glyphs. So any text used for searching /
# length-measuring against a value obtained from ".Text" must use plain
# ASCII quotes; the curly glyphs are only needed in the strings we *write*.

# =========================================================================
# Slide 2 ("Questions:") - content placeholder
#   Paragraph 1 reworded & split into 3 runs; Paragraph 2 reworded (1 run)
# =========================================================================
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

# --- Paragraph 1: rewrite full paragraph text, then split into 3 runs ---
$oldPara1 = "Why the coyote couldn't memorize the ""song"" although he tried so hard?"
$newPara1 = "Why couldn’t the coyote memorize the “song” although he tried so hard?"
$whole1 = $tr2.Characters(1, $oldPara1.Length)
$whole1.Text = $newPara1

$run1a = "Why "
$run1b = "couldn’t the "
$run1c = "coyote memorize the “song” although he tried so hard?"

$c1 = $tr2.Characters(1, $run1a.Length)
$c1.Text = $run1a

$c2 = $tr2.Characters($run1a.Length + 1, $run1b.Length)
$c2.Text = $run1b

$c3 = $tr2.Characters($run1a.Length + $run1b.Length + 1, $run1c.Length)
$c3.Text = $run1c

# --- Paragraph 2: rewrite full paragraph text (single run) ---
$oldPara2 = "Why the coyote misunderstand the true meaning of his song?"
$newPara2 = "Why did the coyote misunderstand the true meaning of the “song”?"
$para2Start = $run1a.Length + $run1b.Length + $run1c.Length + 2
$whole2 = $tr2.Characters($para2Start, $oldPara2.Length)
$whole2.Text = $newPara2

# =========================================================================
# Slide 3 - content placeholder, paragraph 3 ("E.g. memorize words ...")
#   "memorize" -> "memorizing" (text tweak, single run stays single run)
# =========================================================================
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

$oldPara3 = "E.g. memorize words in contents; understanding the meaning of a poem before reciting it."
$newPara3 = "E.g. memorizing words in contents; understanding the meaning of a poem before reciting it."

$full3 = $tr3.Text
$start3 = $full3.IndexOf($oldPara3) + 1
$whole3 = $tr3.Characters($start3, $oldPara3.Length)
$whole3.Text = $newPara3

# =========================================================================
# Slide 5 - content placeholder, paragraph 3 ("If you don't, then ...")
#   Split the trailing "Example, ..." clause off into its own run and
#   shorten "Example," to "E.g.,"
# =========================================================================
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange

$run3a = "If you don’t, then you should actively get involved. "
$oldTail = "Example, when learning physics, doing some experiments to verify the laws is extremely helpful."
$run3b = "E.g., "
$run3c = "when learning physics, doing some experiments to verify the laws is extremely helpful."

# plain-ASCII version of run3a, for matching against ".Text" output only
$run3aPlain = "If you don't, then you should actively get involved. "

$full5 = $tr5.Text
$start5 = $full5.IndexOf($run3aPlain + $oldTail) + 1

# Rewrite the whole paragraph text first (Example -> E.g.,)
$wholeLen5 = $run3aPlain.Length + $oldTail.Length
$whole5 = $tr5.Characters($start5, $wholeLen5)
$whole5.Text = $run3a + $run3b + $run3c

# Now split it into 3 runs
$d1 = $tr5.Characters($start5, $run3a.Length)
$d1.Text = $run3a

$d2 = $tr5.Characters($start5 + $run3a.Length, $run3b.Length)
$d2.Text = $run3b

$d3 = $tr5.Characters($start5 + $run3a.Length + $run3b.Length, $run3c.Length)
$d3.Text = $run3c

Write-Host "Slide2 final: [$($tr2.Text)]"
Write-Host "Slide3 final: [$($tr3.Text)]"
Write-Host "Slide5 final: [$($tr5.Text)]"
